# Auto-generated edit script applying the Maduin_Profits market-price refresh
# (H/I/J/K/L/M/N recalculated columns) described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 942.5
$ws.Range("J17").Value = 942.5
$ws.Range("L17").Value = 2827.5
$ws.Range("N17").Value = -3163.5
$ws.Range("H33").Value = 444.8125
$ws.Range("I33").Value = 424.42856
$ws.Range("K33").Value = 424.42856
$ws.Range("M33").Value = -195.42856
$ws.Range("H62").Value = 9722
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376
$ws.Range("H65").Value = 9722
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880
$ws.Range("H86").Value = 2633.2942
$ws.Range("I86").Value = 3655.25
$ws.Range("J86").Value = 1724.8889
$ws.Range("K86").Value = 3655.25
$ws.Range("L86").Value = 1724.8889
$ws.Range("M86").Value = -2532.25
$ws.Range("N86").Value = -3970.8889
$ws.Range("H88").Value = 4082.889
$ws.Range("J88").Value = 4093.5
$ws.Range("L88").Value = 4093.5
$ws.Range("N88").Value = -4905.5
$ws.Range("H89").Value = 2633.2942
$ws.Range("I89").Value = 3655.25
$ws.Range("J89").Value = 1724.8889
$ws.Range("K89").Value = 18276.25
$ws.Range("L89").Value = 8624.4445
$ws.Range("M89").Value = -12660.25
$ws.Range("N89").Value = -19856.4445
$ws.Range("H91").Value = 4082.889
$ws.Range("J91").Value = 4093.5
$ws.Range("L91").Value = 4093.5
$ws.Range("N91").Value = -6901.5
$ws.Range("H96").Value = 1238.1428
$ws.Range("I96").Value = 644.75
$ws.Range("J96").Value = 2029.3334
$ws.Range("K96").Value = 1934.25
$ws.Range("L96").Value = 6088.0002
$ws.Range("M96").Value = -561.25
$ws.Range("N96").Value = -8834.0002
$ws.Range("H132").Value = 4612.227
$ws.Range("I132").Value = 3820.5386
$ws.Range("K132").Value = 11461.6158
$ws.Range("M132").Value = -8931.6158
$ws.Range("H137").Value = 2562.5
$ws.Range("I137").Value = 2300
$ws.Range("K137").Value = 6900
$ws.Range("M137").Value = -4350
$ws.Range("H138").Value = 3139.4614
$ws.Range("I138").Value = 2801.2727
$ws.Range("K138").Value = 8403.8181
$ws.Range("M138").Value = -3263.8181
$ws.Range("H141").Value = 1551.5385
$ws.Range("I141").Value = 1597.5
$ws.Range("K141").Value = 4792.5
$ws.Range("M141").Value = 387.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9505.5
$ws.Range("I2").Value = 9000
$ws.Range("J2").Value = 10011
$ws.Range("K2").Value = 9000
$ws.Range("L2").Value = 10011
$ws.Range("M2").Value = -8887
$ws.Range("N2").Value = -10237
$ws.Range("H32").Value = 3067.973
$ws.Range("I32").Value = 2120.875
$ws.Range("K32").Value = 2120.875
$ws.Range("M32").Value = -1833.875
$ws.Range("H63").Value = 4777.8
$ws.Range("I63").Value = 2696
$ws.Range("J63").Value = 6165.6665
$ws.Range("K63").Value = 2696
$ws.Range("L63").Value = 6165.6665
$ws.Range("M63").Value = -2010
$ws.Range("N63").Value = -7537.6665
$ws.Range("H66").Value = 4777.8
$ws.Range("I66").Value = 2696
$ws.Range("J66").Value = 6165.6665
$ws.Range("K66").Value = 13480
$ws.Range("L66").Value = 30828.3325
$ws.Range("M66").Value = -10048
$ws.Range("N66").Value = -37692.3325
$ws.Range("H88").Value = 5556.8887
$ws.Range("J88").Value = 5858.143
$ws.Range("L88").Value = 5858.143
$ws.Range("N88").Value = -6670.143
$ws.Range("H91").Value = 5556.8887
$ws.Range("J91").Value = 5858.143
$ws.Range("L91").Value = 5858.143
$ws.Range("N91").Value = -8666.143
$ws.Range("H97").Value = 1038.8948
$ws.Range("I97").Value = 850.6
$ws.Range("J97").Value = 1745
$ws.Range("K97").Value = 850.6
$ws.Range("L97").Value = 1745
$ws.Range("M97").Value = -354.6
$ws.Range("N97").Value = -2737
$ws.Range("H103").Value = 15000
$ws.Range("J103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("N103").Value = -17344
$ws.Range("H110").Value = 1032
$ws.Range("I110").Value = 1032
$ws.Range("K110").Value = 1032
$ws.Range("M110").Value = 1013
$ws.Range("H116").Value = 9505.5
$ws.Range("I116").Value = 9000
$ws.Range("J116").Value = 10011
$ws.Range("K116").Value = 9000
$ws.Range("L116").Value = 10011
$ws.Range("M116").Value = -6706
$ws.Range("N116").Value = -14599
$ws.Range("H119").Value = 52200
$ws.Range("J119").Value = 52200
$ws.Range("L119").Value = 52200
$ws.Range("N119").Value = -61876
$ws.Range("H122").Value = 5506.5
$ws.Range("I122").Value = 4999
$ws.Range("J122").Value = 6014
$ws.Range("K122").Value = 14997
$ws.Range("L122").Value = 18042
$ws.Range("M122").Value = -12547
$ws.Range("N122").Value = -22942

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9505.5
$ws.Range("I3").Value = 9000
$ws.Range("J3").Value = 10011
$ws.Range("K3").Value = 9000
$ws.Range("L3").Value = 10011
$ws.Range("M3").Value = -8886
$ws.Range("N3").Value = -10239
$ws.Range("H22").Value = 660
$ws.Range("I22").Value = 660
$ws.Range("K22").Value = 660
$ws.Range("M22").Value = -487
$ws.Range("H74").Value = 53589.668
$ws.Range("J74").Value = 53589.668
$ws.Range("L74").Value = 53589.668
$ws.Range("N74").Value = -55461.668
$ws.Range("H77").Value = 53589.668
$ws.Range("J77").Value = 53589.668
$ws.Range("L77").Value = 160769.004
$ws.Range("N77").Value = -170129.004
$ws.Range("H94").Value = 3450.1667
$ws.Range("I94").Value = 2854.077
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 2854.077
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -2403.077
$ws.Range("N94").Value = -5902
$ws.Range("H107").Value = 944.875
$ws.Range("I107").Value = 866.8
$ws.Range("J107").Value = 1075
$ws.Range("K107").Value = 866.8
$ws.Range("L107").Value = 1075
$ws.Range("M107").Value = 1053.2
$ws.Range("N107").Value = -4915

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 410
$ws.Range("J68").Value = 410
$ws.Range("L68").Value = 1230
$ws.Range("N68").Value = -2852
$ws.Range("H71").Value = 410
$ws.Range("J71").Value = 410
$ws.Range("L71").Value = 3690
$ws.Range("N71").Value = -11802
$ws.Range("H107").Value = 745.6667
$ws.Range("J107").Value = 766.5454999999999
$ws.Range("L107").Value = 2299.6365
$ws.Range("N107").Value = -6139.6365
$ws.Range("H122").Value = 492.14285
$ws.Range("I122").Value = 450.2857
$ws.Range("J122").Value = 534
$ws.Range("K122").Value = 4052.5713
$ws.Range("L122").Value = 4806
$ws.Range("M122").Value = -1602.5713
$ws.Range("N122").Value = -9706

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3049.8
$ws.Range("I80").Value = 2799
$ws.Range("J80").Value = 3112.5
$ws.Range("K80").Value = 2799
$ws.Range("L80").Value = 3112.5
$ws.Range("M80").Value = -1801
$ws.Range("N80").Value = -5108.5
$ws.Range("H83").Value = 3049.8
$ws.Range("I83").Value = 2799
$ws.Range("J83").Value = 3112.5
$ws.Range("K83").Value = 13995
$ws.Range("L83").Value = 15562.5
$ws.Range("M83").Value = -9003
$ws.Range("N83").Value = -25546.5
$ws.Range("H113").Value = 2220
$ws.Range("I113").Value = 2293.3333
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2293.3333
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -123.3332999999998
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 1499.75
$ws.Range("I122").Value = 1499.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4499.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2049.25
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1298.6
$ws.Range("I7").Value = 1348.25
$ws.Range("K7").Value = 1348.25
$ws.Range("M7").Value = -1236.25
$ws.Range("H22").Value = 1799.7778
$ws.Range("I22").Value = 1024.75
$ws.Range("K22").Value = 1024.75
$ws.Range("M22").Value = -729.75
$ws.Range("H27").Value = 1799.7778
$ws.Range("I27").Value = 1024.75
$ws.Range("K27").Value = 1024.75
$ws.Range("M27").Value = -917.75
$ws.Range("H46").Value = 1753.3
$ws.Range("I46").Value = 1294.6842
$ws.Range("K46").Value = 1294.6842
$ws.Range("M46").Value = -1106.6842
$ws.Range("H122").Value = 5661.5
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H126").Value = 1298.6
$ws.Range("I126").Value = 1348.25
$ws.Range("K126").Value = 4044.75
$ws.Range("M126").Value = -1574.75
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H96").Value = 2766.6667
$ws.Range("I96").Value = 2200
$ws.Range("J96").Value = 3333.3333
$ws.Range("K96").Value = 2200
$ws.Range("L96").Value = 3333.3333
$ws.Range("M96").Value = -827
$ws.Range("N96").Value = -6079.3333
$ws.Range("H119").Value = 49749.25
$ws.Range("J119").Value = 49749.25
$ws.Range("L119").Value = 49749.25
$ws.Range("N119").Value = -59425.25
$ws.Range("H122").Value = 4666.3335
$ws.Range("I122").Value = 4499.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 13498.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -11048.5
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 1610.2
$ws.Range("I132").Value = 675.5
$ws.Range("K132").Value = 2026.5
$ws.Range("M132").Value = 503.5
